# Week 3 Stat Corrections
# - Anthony's (row 4) Week 3 score corrected from 85 to 110 (Total 110 -> 135)
# - Dean's (row 10) Week 3 score corrected from 45 to 55 (Total 130 -> 140)
# All downstream VLOOKUP/SUM cells on the league sheets recalc automatically.

$wb = $excel.ActiveWorkbook

$wsPoints = $wb.Worksheets.Item("Contestant Points")
$wsEdX    = $wb.Worksheets.Item("edX League")

# --- Week 3 ("S" column) score corrections -------------------------------
$wsPoints.Range("S4").Value  = 110   # Anthony
$wsPoints.Range("S10").Value = 55    # Dean

# --- Re-enter the Total Points formula across the whole column so Excel
#     collapses it into a single shared formula (matches the saved file's
#     <f t="shared" ref="T2:T32" si="0">SUM(Q2:S2)</f> pattern). -----------
$wsPoints.Range("T2:T32").Formula = "=SUM(Q2:S2)"

# --- Column width touch-up on "Contestant Points" -------------------------
# Columns B:G end up as a single uniform custom-width run (9.5546875 chars).
$wsPoints.Range("B1:G1").ColumnWidth = 9.5546875

# --- Sheet/selection/active-tab bookkeeping -------------------------------
# The author ended the session with "Contestant Points" focused (cell G10
# selected) instead of "edX League" (which now just has C8 selected).
# Set the no-longer-active sheet's selection first, then finish on the
# sheet that should end up active/selected (selecting a range activates
# its sheet as a side effect).
$wsEdX.Range("C8").Select()

$wsPoints.Activate()
$wsPoints.Range("G10").Select()
